$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: SCD0170 -> SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID cell (B2): DGS-185 -> SCD0011-001
$ws.Range("B2").Value = "SCD0011-001"

# Reflect the saved cursor/selection state: active cell moved to B3
# (also resets the scrolled "topLeftCell" back to default)
$ws.Range("B3").Select()
